$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.206015333333333
$ws.Range("H2").Value = 18.618046
$ws.Range("I2").Value = 0.0150172404156507
$ws.Range("J2").Value = 0.0150172404156507
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1575256666666667
$ws.Range("N2").Value = 0.472577
$ws.Range("O2").Value = 0.6985926944284299
$ws.Range("P2").Value = 0.69859269442843
$ws.Range("Q2").Value = 0.9776067027268889
$ws.Range("R2").Value = 8.798460324542001
$ws.Range("S2").Value = 0.01049093444484894
$ws.Range("T2").Value = 0.01049093444484894

$ws.Range("G3").Value = 6.206015333333333
$ws.Range("H3").Value = 18.618046
$ws.Range("I3").Value = 0.0150172404156507
$ws.Range("J3").Value = 0.0150172404156507
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.06796433333333333
$ws.Range("N3").Value = 0.203893
$ws.Range("O3").Value = 0.30140730557157
$ws.Range("P3").Value = 0.30140730557157
$ws.Range("Q3").Value = 0.4217876947864445
$ws.Range("R3").Value = 3.796089253078
$ws.Range("S3").Value = 0.004526305970801761
$ws.Range("T3").Value = 0.004526305970801761

$ws.Range("I4").Value = 0.9317452840597572
$ws.Range("J4").Value = 0.9317452840597571
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1575256666666667
$ws.Range("N4").Value = 0.472577
$ws.Range("O4").Value = 0.6985926944284299
$ws.Range("P4").Value = 0.69859269442843
$ws.Range("Q4").Value = 60.65564709090523
$ws.Range("R4").Value = 545.900823818147
$ws.Range("S4").Value = 0.6509104485122885
$ws.Range("T4").Value = 0.6509104485122885

$ws.Range("I5").Value = 0.9317452840597572
$ws.Range("J5").Value = 0.9317452840597571
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06796433333333333
$ws.Range("N5").Value = 0.203893
$ws.Range("O5").Value = 0.30140730557157
$ws.Range("P5").Value = 0.30140730557157
$ws.Range("Q5").Value = 26.16983444455811
$ws.Range("R5").Value = 235.528510001023
$ws.Range("S5").Value = 0.2808348355474686
$ws.Range("T5").Value = 0.2808348355474685

$ws.Range("G6").Value = 22.00088566666667
$ws.Range("H6").Value = 66.002657
$ws.Range("I6").Value = 0.05323747552459213
$ws.Range("J6").Value = 0.05323747552459213
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1575256666666667
$ws.Range("N6").Value = 0.472577
$ws.Range("O6").Value = 0.6985926944284299
$ws.Range("P6").Value = 0.69859269442843
$ws.Range("Q6").Value = 3.465704181898778
$ws.Range("R6").Value = 31.191337637089
$ws.Range("S6").Value = 0.03719131147129241
$ws.Range("T6").Value = 0.03719131147129241

$ws.Range("G7").Value = 22.00088566666667
$ws.Range("H7").Value = 66.002657
$ws.Range("I7").Value = 0.05323747552459213
$ws.Range("J7").Value = 0.05323747552459213
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.06796433333333333
$ws.Range("N7").Value = 0.203893
$ws.Range("O7").Value = 0.30140730557157
$ws.Range("P7").Value = 0.30140730557157
$ws.Range("Q7").Value = 1.495275527077889
$ws.Range("R7").Value = 13.457479743701
$ws.Range("S7").Value = 0.01604616405329972
$ws.Range("T7").Value = 0.01604616405329972
